$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Population Definitions")
Write-Host $ws.Name
